$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column F header, matching the style used by the other header cells (e.g. E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the time_taken values for each data row
$ws.Range("F2").Value = "2021-10-05 13:38:32.782107"
$ws.Range("F3").Value = "2021-10-05 13:38:32.782119"
$ws.Range("F4").Value = "2021-10-05 13:38:32.782123"
$ws.Range("F5").Value = "2021-10-05 13:38:32.782126"

$wb.Save()
